$wb = $excel.ActiveWorkbook

# --- BBNPPTY sheet: set B,C,D (years 2021-2023) to 0 for all fuel types ---
# except "hard coal w CCS", "natural gas combined cycle w CCS", "biomass w CCS",
# "lignite w CCS" (data rows 19-22), which remain 1.
$wsData = $wb.Worksheets.Item("BBNPPTY")

$rowsToZero = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,23,24,25)
foreach ($r in $rowsToZero) {
    $wsData.Range("B" + $r + ":D" + $r).Value = 0
}

# --- View / selection changes ---
# "About" sheet loses tabSelected, selection unaffected otherwise.
$wsAbout = $wb.Worksheets.Item("About")
$null = $wsAbout.Range("B22").Select()

# "BBNPPTY" becomes the active/selected tab, with a new selection.
$wsData.Activate()
$null = $wsData.Range("G27").Select()
